$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(1, 1).Value = 'Marjorie Mann'
$ws.Cells.Item(2, 1).Value = 'Mrs. Kathryne Shanahan IV'
$ws.Cells.Item(3, 1).Value = 'Isabel Satterfield'
$ws.Cells.Item(4, 1).Value = 'Junius Christiansen'
$ws.Cells.Item(5, 1).Value = 'Mr. Lon Marquardt Jr.'
$ws.Cells.Item(6, 1).Value = 'Virgie Feest'
$ws.Cells.Item(7, 1).Value = 'Roderick Eichmann'
$ws.Cells.Item(8, 1).Value = 'Raleigh Lynch'
$ws.Cells.Item(9, 1).Value = 'Prof. Jerry Murray'
$ws.Cells.Item(10, 1).Value = 'Guy Farrell'
$ws.Cells.Item(11, 1).Value = 'Dr. Edwardo Heaney III'
$ws.Cells.Item(12, 1).Value = 'Valentine Kub'
$ws.Cells.Item(13, 1).Value = 'Zola Lehner'
$ws.Cells.Item(14, 1).Value = 'Mark Marvin'
$ws.Cells.Item(15, 1).Value = 'Earl Tillman Sr.'
$ws.Cells.Item(16, 1).Value = 'Ottilie Treutel'
$ws.Cells.Item(17, 1).Value = 'Kenton Dickens'
$ws.Cells.Item(18, 1).Value = 'Moises Wunsch'
$ws.Cells.Item(19, 1).Value = 'Dr. Ena Pouros PhD'
$ws.Cells.Item(20, 1).Value = 'Miss Adrienne Olson'
$ws.Cells.Item(21, 1).Value = 'Prof. Kieran Roob Jr.'
$ws.Cells.Item(22, 1).Value = 'Agustin Gutmann'
$ws.Cells.Item(23, 1).Value = 'Mrs. Marisol O''Conner'
$ws.Cells.Item(24, 1).Value = 'Parker Kreiger'
$ws.Cells.Item(25, 1).Value = 'Linnie Stoltenberg Jr.'
$ws.Cells.Item(26, 1).Value = 'Bulah Davis'
$ws.Cells.Item(27, 1).Value = 'Nathanael Cole'
$ws.Cells.Item(28, 1).Value = 'Emie Thiel'
$ws.Cells.Item(29, 1).Value = 'Jolie Eichmann'
$ws.Cells.Item(30, 1).Value = 'Miss Ofelia Bailey IV'
$ws.Cells.Item(31, 1).Value = 'Maxime Langworth'
$ws.Cells.Item(32, 1).Value = 'Jerrold Gutkowski'
$ws.Cells.Item(33, 1).Value = 'Dr. Kyla Haag'
$ws.Cells.Item(34, 1).Value = 'Lola Becker DVM'
$ws.Cells.Item(35, 1).Value = 'Lura Kuhic'
$ws.Cells.Item(36, 1).Value = 'Wendell Purdy'
$ws.Cells.Item(37, 1).Value = 'Meta Ruecker'
$ws.Cells.Item(38, 1).Value = 'Mr. Carson Flatley IV'
$ws.Cells.Item(39, 1).Value = 'Jakob Gorczany'
$ws.Cells.Item(40, 1).Value = 'Ms. Blanca Mueller'
$ws.Cells.Item(41, 1).Value = 'Damon Hoeger PhD'
$ws.Cells.Item(42, 1).Value = 'Leonie Mueller'
$ws.Cells.Item(43, 1).Value = 'Colby Koss'
$ws.Cells.Item(44, 1).Value = 'Samir Fadel'
$ws.Cells.Item(45, 1).Value = 'Prof. Jevon Sporer'
$ws.Cells.Item(46, 1).Value = 'Verner Gutkowski'
$ws.Cells.Item(47, 1).Value = 'Lillian Christiansen III'
$ws.Cells.Item(48, 1).Value = 'Lexie Morissette'
$ws.Cells.Item(49, 1).Value = 'Leone Miller'
$ws.Cells.Item(50, 1).Value = 'Gideon Klocko'
$ws.Cells.Item(51, 1).Value = 'Demario Waelchi'
$ws.Cells.Item(52, 1).Value = 'Jamel Gottlieb'
$ws.Cells.Item(53, 1).Value = 'Dayana Mosciski'
$ws.Cells.Item(54, 1).Value = 'Harley Donnelly DDS'
$ws.Cells.Item(55, 1).Value = 'Allie Hill'
$ws.Cells.Item(56, 1).Value = 'Miss Simone Larson IV'
$ws.Cells.Item(57, 1).Value = 'Kirk Kutch'
$ws.Cells.Item(58, 1).Value = 'Mr. Frederick Howe'
$ws.Cells.Item(59, 1).Value = 'Gwen Toy IV'
$ws.Cells.Item(60, 1).Value = 'Ms. Mabelle Hettinger'
$ws.Cells.Item(61, 1).Value = 'Simeon Bosco'
$ws.Cells.Item(62, 1).Value = 'Duncan Frami IV'
$ws.Cells.Item(63, 1).Value = 'Miss Carmella King DVM'
$ws.Cells.Item(64, 1).Value = 'Mohammad Wolff'
$ws.Cells.Item(65, 1).Value = 'Mr. Antone Beahan Jr.'
$ws.Cells.Item(66, 1).Value = 'Mrs. Palma Stracke'
$ws.Cells.Item(67, 1).Value = 'Malvina Stehr'
$ws.Cells.Item(68, 1).Value = 'Miss Abigale Corkery'
$ws.Cells.Item(69, 1).Value = 'Floyd Kub'
$ws.Cells.Item(70, 1).Value = 'Ms. Delpha Murphy'
$ws.Cells.Item(71, 1).Value = 'Dillan Wyman'
$ws.Cells.Item(72, 1).Value = 'Dr. Ole Gleichner DVM'
$ws.Cells.Item(73, 1).Value = 'Mrs. Daniela Gleason'
$ws.Cells.Item(74, 1).Value = 'Fredy Walsh'
$ws.Cells.Item(75, 1).Value = 'Laila Considine'
$ws.Cells.Item(76, 1).Value = 'Triston Bauch'
$ws.Cells.Item(77, 1).Value = 'Joel Turner II'
$ws.Cells.Item(78, 1).Value = 'Oswald Lubowitz Sr.'
$ws.Cells.Item(79, 1).Value = 'Margarete Hill'
$ws.Cells.Item(80, 1).Value = 'Emmy Bahringer'
$ws.Cells.Item(81, 1).Value = 'Prof. Pansy Johnson IV'
$ws.Cells.Item(82, 1).Value = 'Jailyn Hackett DVM'
$ws.Cells.Item(83, 1).Value = 'Phoebe Walker'
$ws.Cells.Item(84, 1).Value = 'Ray Corkery'
$ws.Cells.Item(85, 1).Value = 'Mr. Waylon Trantow'
$ws.Cells.Item(86, 1).Value = 'Prof. Dorian Barrows'
$ws.Cells.Item(87, 1).Value = 'Emmanuel Kris PhD'
$ws.Cells.Item(88, 1).Value = 'Lenna Prohaska'
$ws.Cells.Item(89, 1).Value = 'Miss Aliya Swaniawski'
$ws.Cells.Item(90, 1).Value = 'Antonio Beatty'
$ws.Cells.Item(91, 1).Value = 'Cooper Becker'
$ws.Cells.Item(92, 1).Value = 'Ben McKenzie DVM'
$ws.Cells.Item(93, 1).Value = 'Dr. Briana Weissnat DVM'
$ws.Cells.Item(94, 1).Value = 'Dr. Trisha McGlynn PhD'
$ws.Cells.Item(95, 1).Value = 'Valerie Osinski'
$ws.Cells.Item(96, 1).Value = 'Keanu Feest'
$ws.Cells.Item(97, 1).Value = 'Fannie Barton'
$ws.Cells.Item(98, 1).Value = 'Yoshiko Fritsch'
$ws.Cells.Item(99, 1).Value = 'Christian Kuhic'
$ws.Cells.Item(100, 1).Value = 'Ms. Carolyn Beer IV'
